# This script applies a row-level permutation to the price-data table
# (rows 2-14, columns D/J/K/L/M/P) to match the weekly data refresh
# described in the commit message ("Fruta / hortaliza, semanal").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture original values (Fecha, Volumen, Precio minimo/maximo/promedio, Precio $/Kg)
# for rows 2-14 before overwriting anything, since the mapping is a permutation.
$orig = @{}
$orig[2] = @{
    "D" = $ws.Range("D2").Value2
    "J" = $ws.Range("J2").Value2
    "K" = $ws.Range("K2").Value2
    "L" = $ws.Range("L2").Value2
    "M" = $ws.Range("M2").Value2
    "P" = $ws.Range("P2").Value2
}
$orig[3] = @{
    "D" = $ws.Range("D3").Value2
    "J" = $ws.Range("J3").Value2
    "K" = $ws.Range("K3").Value2
    "L" = $ws.Range("L3").Value2
    "M" = $ws.Range("M3").Value2
    "P" = $ws.Range("P3").Value2
}
$orig[4] = @{
    "D" = $ws.Range("D4").Value2
    "J" = $ws.Range("J4").Value2
    "K" = $ws.Range("K4").Value2
    "L" = $ws.Range("L4").Value2
    "M" = $ws.Range("M4").Value2
    "P" = $ws.Range("P4").Value2
}
$orig[5] = @{
    "D" = $ws.Range("D5").Value2
    "J" = $ws.Range("J5").Value2
    "K" = $ws.Range("K5").Value2
    "L" = $ws.Range("L5").Value2
    "M" = $ws.Range("M5").Value2
    "P" = $ws.Range("P5").Value2
}
$orig[6] = @{
    "D" = $ws.Range("D6").Value2
    "J" = $ws.Range("J6").Value2
    "K" = $ws.Range("K6").Value2
    "L" = $ws.Range("L6").Value2
    "M" = $ws.Range("M6").Value2
    "P" = $ws.Range("P6").Value2
}
$orig[7] = @{
    "D" = $ws.Range("D7").Value2
    "J" = $ws.Range("J7").Value2
    "K" = $ws.Range("K7").Value2
    "L" = $ws.Range("L7").Value2
    "M" = $ws.Range("M7").Value2
    "P" = $ws.Range("P7").Value2
}
$orig[8] = @{
    "D" = $ws.Range("D8").Value2
    "J" = $ws.Range("J8").Value2
    "K" = $ws.Range("K8").Value2
    "L" = $ws.Range("L8").Value2
    "M" = $ws.Range("M8").Value2
    "P" = $ws.Range("P8").Value2
}
$orig[9] = @{
    "D" = $ws.Range("D9").Value2
    "J" = $ws.Range("J9").Value2
    "K" = $ws.Range("K9").Value2
    "L" = $ws.Range("L9").Value2
    "M" = $ws.Range("M9").Value2
    "P" = $ws.Range("P9").Value2
}
$orig[10] = @{
    "D" = $ws.Range("D10").Value2
    "J" = $ws.Range("J10").Value2
    "K" = $ws.Range("K10").Value2
    "L" = $ws.Range("L10").Value2
    "M" = $ws.Range("M10").Value2
    "P" = $ws.Range("P10").Value2
}
$orig[11] = @{
    "D" = $ws.Range("D11").Value2
    "J" = $ws.Range("J11").Value2
    "K" = $ws.Range("K11").Value2
    "L" = $ws.Range("L11").Value2
    "M" = $ws.Range("M11").Value2
    "P" = $ws.Range("P11").Value2
}
$orig[12] = @{
    "D" = $ws.Range("D12").Value2
    "J" = $ws.Range("J12").Value2
    "K" = $ws.Range("K12").Value2
    "L" = $ws.Range("L12").Value2
    "M" = $ws.Range("M12").Value2
    "P" = $ws.Range("P12").Value2
}
$orig[13] = @{
    "D" = $ws.Range("D13").Value2
    "J" = $ws.Range("J13").Value2
    "K" = $ws.Range("K13").Value2
    "L" = $ws.Range("L13").Value2
    "M" = $ws.Range("M13").Value2
    "P" = $ws.Range("P13").Value2
}
$orig[14] = @{
    "D" = $ws.Range("D14").Value2
    "J" = $ws.Range("J14").Value2
    "K" = $ws.Range("K14").Value2
    "L" = $ws.Range("L14").Value2
    "M" = $ws.Range("M14").Value2
    "P" = $ws.Range("P14").Value2
}

# Write back the values following the new row order.
# Row 2 now holds the data previously found in row 8
$ws.Range("D2").Value = $orig[8]["D"]
$ws.Range("J2").Value = $orig[8]["J"]
$ws.Range("K2").Value = $orig[8]["K"]
$ws.Range("L2").Value = $orig[8]["L"]
$ws.Range("M2").Value = $orig[8]["M"]
$ws.Range("P2").Value = $orig[8]["P"]

# Row 3 now holds the data previously found in row 5
$ws.Range("D3").Value = $orig[5]["D"]
$ws.Range("J3").Value = $orig[5]["J"]
$ws.Range("K3").Value = $orig[5]["K"]
$ws.Range("L3").Value = $orig[5]["L"]
$ws.Range("M3").Value = $orig[5]["M"]
$ws.Range("P3").Value = $orig[5]["P"]

# Row 4 now holds the data previously found in row 7
$ws.Range("D4").Value = $orig[7]["D"]
$ws.Range("J4").Value = $orig[7]["J"]
$ws.Range("K4").Value = $orig[7]["K"]
$ws.Range("L4").Value = $orig[7]["L"]
$ws.Range("M4").Value = $orig[7]["M"]
$ws.Range("P4").Value = $orig[7]["P"]

# Row 5 now holds the data previously found in row 14
$ws.Range("D5").Value = $orig[14]["D"]
$ws.Range("J5").Value = $orig[14]["J"]
$ws.Range("K5").Value = $orig[14]["K"]
$ws.Range("L5").Value = $orig[14]["L"]
$ws.Range("M5").Value = $orig[14]["M"]
$ws.Range("P5").Value = $orig[14]["P"]

# Row 6 now holds the data previously found in row 9
$ws.Range("D6").Value = $orig[9]["D"]
$ws.Range("J6").Value = $orig[9]["J"]
$ws.Range("K6").Value = $orig[9]["K"]
$ws.Range("L6").Value = $orig[9]["L"]
$ws.Range("M6").Value = $orig[9]["M"]
$ws.Range("P6").Value = $orig[9]["P"]

# Row 7 now holds the data previously found in row 6
$ws.Range("D7").Value = $orig[6]["D"]
$ws.Range("J7").Value = $orig[6]["J"]
$ws.Range("K7").Value = $orig[6]["K"]
$ws.Range("L7").Value = $orig[6]["L"]
$ws.Range("M7").Value = $orig[6]["M"]
$ws.Range("P7").Value = $orig[6]["P"]

# Row 8 now holds the data previously found in row 4
$ws.Range("D8").Value = $orig[4]["D"]
$ws.Range("J8").Value = $orig[4]["J"]
$ws.Range("K8").Value = $orig[4]["K"]
$ws.Range("L8").Value = $orig[4]["L"]
$ws.Range("M8").Value = $orig[4]["M"]
$ws.Range("P8").Value = $orig[4]["P"]

# Row 9 now holds the data previously found in row 3
$ws.Range("D9").Value = $orig[3]["D"]
$ws.Range("J9").Value = $orig[3]["J"]
$ws.Range("K9").Value = $orig[3]["K"]
$ws.Range("L9").Value = $orig[3]["L"]
$ws.Range("M9").Value = $orig[3]["M"]
$ws.Range("P9").Value = $orig[3]["P"]

# Row 10 now holds the data previously found in row 2
$ws.Range("D10").Value = $orig[2]["D"]
$ws.Range("J10").Value = $orig[2]["J"]
$ws.Range("K10").Value = $orig[2]["K"]
$ws.Range("L10").Value = $orig[2]["L"]
$ws.Range("M10").Value = $orig[2]["M"]
$ws.Range("P10").Value = $orig[2]["P"]

# Row 11 now holds the data previously found in row 12
$ws.Range("D11").Value = $orig[12]["D"]
$ws.Range("J11").Value = $orig[12]["J"]
$ws.Range("K11").Value = $orig[12]["K"]
$ws.Range("L11").Value = $orig[12]["L"]
$ws.Range("M11").Value = $orig[12]["M"]
$ws.Range("P11").Value = $orig[12]["P"]

# Row 12 now holds the data previously found in row 13
$ws.Range("D12").Value = $orig[13]["D"]
$ws.Range("J12").Value = $orig[13]["J"]
$ws.Range("K12").Value = $orig[13]["K"]
$ws.Range("L12").Value = $orig[13]["L"]
$ws.Range("M12").Value = $orig[13]["M"]
$ws.Range("P12").Value = $orig[13]["P"]

# Row 13 now holds the data previously found in row 11
$ws.Range("D13").Value = $orig[11]["D"]
$ws.Range("J13").Value = $orig[11]["J"]
$ws.Range("K13").Value = $orig[11]["K"]
$ws.Range("L13").Value = $orig[11]["L"]
$ws.Range("M13").Value = $orig[11]["M"]
$ws.Range("P13").Value = $orig[11]["P"]

# Row 14 now holds the data previously found in row 10
$ws.Range("D14").Value = $orig[10]["D"]
$ws.Range("J14").Value = $orig[10]["J"]
$ws.Range("K14").Value = $orig[10]["K"]
$ws.Range("L14").Value = $orig[10]["L"]
$ws.Range("M14").Value = $orig[10]["M"]
$ws.Range("P14").Value = $orig[10]["P"]

